$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "Categories" column (J) with the group/keyword category labels
# produced by handle_obvious_dups().
$ws.Range("J2").Value = "1 catA"
$ws.Range("J3").Value = "1 catA"
$ws.Range("J4").Value = "2 catB"
$ws.Range("J5").Value = "1 catA"
$ws.Range("J6").Value = "2 catB"

$ws.Range("A1").Select() | Out-Null
